$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank separator row (row 3) - clearing it entirely so it
# disappears from the sheet (row numbers below are NOT shifted).
$ws.Range("A3:B3").Clear()

# Center-align the serial number column for the already numbered rows.
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A4:A22").HorizontalAlignment = -4108

# Row 23 (Slno 21) gets bold + centered Slno, and bold school name.
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").HorizontalAlignment = -4108
$ws.Range("B23").Font.Bold = $true

# Rows 24-30 previously had no serial numbers; add them (22-28),
# matching the bold/centered style used for the rest of the list, and
# make the school name bold to match the other entries.
for ($row = 24; $row -le 30; $row++) {
    $slno = $row - 2
    $ws.Cells.Item($row, 1).Value = $slno
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 2).Font.Bold = $true
}

# Update the view: select B31 (the next empty row), as left by the editor.
$ws.Range("B31").Select()
